# Apply the "Update countries & provincias Spain" data refresh to the
# Pais worksheet: country case counts moved around in rank (because the
# underlying totals changed), so some rows need a new country name in
# column A plus refreshed stats in B:H, while other rows only get
# refreshed stats. The "last updated" timestamp cell also changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1 -> Datos actualizados a 1 de Abril de 2020 a las 11:20
$ws.Range("A1").Value = 'Datos actualizados a 1 de Abril de 2020 a las 11:20'

# Row 13 -> Belgica
$ws.Range("A13").Value = 'Belgica'
$ws.Range("B13").Value = 13964
$ws.Range("C13").Value = 1189
$ws.Range("D13").Value = 2132
$ws.Range("E13").Value = 11004
$ws.Range("F13").Value = 1088
$ws.Range("G13").Value = 123
$ws.Range("H13").Value = 828

# Row 14 -> Turquia
$ws.Range("A14").Value = 'Turquia'
$ws.Range("B14").Value = 13531
$ws.Range("D14").Value = 243
$ws.Range("E14").Value = 13074
$ws.Range("F14").Value = 847
$ws.Range("H14").Value = 214

# Row 27 -> Malasia
$ws.Range("A27").Value = 'Malasia'
$ws.Range("B27").Value = 2908
$ws.Range("C27").Value = 142
$ws.Range("D27").Value = 645
$ws.Range("E27").Value = 2218
$ws.Range("F27").Value = 102
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 45

# Row 28 -> Dinamarca
$ws.Range("A28").Value = 'Dinamarca'
$ws.Range("B28").Value = 2860
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 2769
$ws.Range("F28").Value = 145
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 90

# Row 29 -> Rusia
$ws.Range("A29").Value = 'Rusia'
$ws.Range("B29").Value = 2777
$ws.Range("C29").Value = 440
$ws.Range("D29").Value = 190
$ws.Range("E29").Value = 2563
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 24

# Row 39 -> Indonesia
$ws.Range("A39").Value = 'Indonesia'
$ws.Range("B39").Value = 1677
$ws.Range("C39").Value = 149
$ws.Range("D39").Value = 103
$ws.Range("E39").Value = 1417
$ws.Range("G39").Value = 21
$ws.Range("H39").Value = 157

# Row 40 -> India
$ws.Range("A40").Value = 'India'
$ws.Range("B40").Value = 1590
$ws.Range("C40").Value = 193
$ws.Range("D40").Value = 148
$ws.Range("E40").Value = 1397
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 45

# Row 41 -> Arabia Saudita
$ws.Range("A41").Value = 'Arabia Saudita'
$ws.Range("B41").Value = 1563
$ws.Range("D41").Value = 165
$ws.Range("E41").Value = 1388
$ws.Range("F41").Value = 31
$ws.Range("H41").Value = 10

# Row 57 (stats only)
$ws.Range("B57").Value = 779
$ws.Range("C57").Value = 34
$ws.Range("E57").Value = 749

# Row 79 (stats only)
$ws.Range("D79").Value = 22
$ws.Range("E79").Value = 327

# Row 117 -> Consejo Danes para los Refugiados
$ws.Range("A117").Value = 'Consejo Danes para los Refugiados'
$ws.Range("C117").Value = 11
$ws.Range("D117").Value = 2
$ws.Range("E117").Value = 99
$ws.Range("F117").Value = 0
$ws.Range("H117").Value = 8

# Row 118 -> Camboya
$ws.Range("A118").Value = 'Camboya'
$ws.Range("B118").Value = 109
$ws.Range("D118").Value = 25
$ws.Range("E118").Value = 84
$ws.Range("F118").Value = 1
$ws.Range("H118").Value = 0

# Row 125 (stats only)
$ws.Range("B125").Value = 65
$ws.Range("C125").Value = 5
$ws.Range("E125").Value = 65

# Row 168 -> Surinam
$ws.Range("A168").Value = 'Surinam'
$ws.Range("C168").Value = 1

# Row 169 -> Seychelles
$ws.Range("A169").Value = 'Seychelles'
$ws.Range("C169").Value = 0

# Row 170 -> Laos
$ws.Range("A170").Value = 'Laos'
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 0
$ws.Range("E170").Value = 10

# Row 171 -> Libia
$ws.Range("A171").Value = 'Libia'
$ws.Range("D171").Value = 1
$ws.Range("E171").Value = 9
$ws.Range("H171").Value = 0

# Row 173 -> Siria
$ws.Range("A173").Value = 'Siria'
$ws.Range("B173").Value = 10
$ws.Range("E173").Value = 8
$ws.Range("H173").Value = 2

# Row 174 -> Granada
$ws.Range("A174").Value = 'Granada'
